$d = $word.ActiveDocument

$d.Content.Find.Execute("54-51=3", $true, $true, $false, $false, $false, $true, 1, $false, "25+68=93", 2) | Out-Null
$d.Content.Find.Execute("19+71=90", $true, $true, $false, $false, $false, $true, 1, $false, "67-19=48", 2) | Out-Null
$d.Content.Find.Execute("72+15=87", $true, $true, $false, $false, $false, $true, 1, $false, "37+29=66", 2) | Out-Null
$d.Content.Find.Execute("8+91=99", $true, $true, $false, $false, $false, $true, 1, $false, "14+32=46", 2) | Out-Null
$d.Content.Find.Execute("56+15=71", $true, $true, $false, $false, $false, $true, 1, $false, "5+76=81", 2) | Out-Null
$d.Content.Find.Execute("99-72=27", $true, $true, $false, $false, $false, $true, 1, $false, "8+54=62", 2) | Out-Null
$d.Content.Find.Execute("48-37=11", $true, $true, $false, $false, $false, $true, 1, $false, "83-45=38", 2) | Out-Null
$d.Content.Find.Execute("12+20=32", $true, $true, $false, $false, $false, $true, 1, $false, "78-46=32", 2) | Out-Null
$d.Content.Find.Execute("6+68=74", $true, $true, $false, $false, $false, $true, 1, $false, "89+4=93", 2) | Out-Null
$d.Content.Find.Execute("79-0=79", $true, $true, $false, $false, $false, $true, 1, $false, "65-16=49", 2) | Out-Null
$d.Content.Find.Execute("93-66=27", $true, $true, $false, $false, $false, $true, 1, $false, "20+11=31", 2) | Out-Null
$d.Content.Find.Execute("78-76=2", $true, $true, $false, $false, $false, $true, 1, $false, "32-16=16", 2) | Out-Null
$d.Content.Find.Execute("17+30=47", $true, $true, $false, $false, $false, $true, 1, $false, "72-23=49", 2) | Out-Null
$d.Content.Find.Execute("65-14=51", $true, $true, $false, $false, $false, $true, 1, $false, "35-4=31", 2) | Out-Null
$d.Content.Find.Execute("24+34=58", $true, $true, $false, $false, $false, $true, 1, $false, "58-32=26", 2) | Out-Null
$d.Content.Find.Execute("24-12=12", $true, $true, $false, $false, $false, $true, 1, $false, "36+27=63", 2) | Out-Null
$d.Content.Find.Execute("55-43=12", $true, $true, $false, $false, $false, $true, 1, $false, "22+52=74", 2) | Out-Null
$d.Content.Find.Execute("49+1=50", $true, $true, $false, $false, $false, $true, 1, $false, "1+72=73", 2) | Out-Null
$d.Content.Find.Execute("27-2=25", $true, $true, $false, $false, $false, $true, 1, $false, "99-90=9", 2) | Out-Null
$d.Content.Find.Execute("65-28=37", $true, $true, $false, $false, $false, $true, 1, $false, "52-15=37", 2) | Out-Null
$d.Content.Find.Execute("40+36=76", $true, $true, $false, $false, $false, $true, 1, $false, "22+3=25", 2) | Out-Null
$d.Content.Find.Execute("87-43=44", $true, $true, $false, $false, $false, $true, 1, $false, "42-6=36", 2) | Out-Null
$d.Content.Find.Execute("1+25=26", $true, $true, $false, $false, $false, $true, 1, $false, "60+5=65", 2) | Out-Null
$d.Content.Find.Execute("5+60=65", $true, $true, $false, $false, $false, $true, 1, $false, "92+1=93", 2) | Out-Null
$d.Content.Find.Execute("69+30=99", $true, $true, $false, $false, $false, $true, 1, $false, "62+3=65", 2) | Out-Null
$d.Content.Find.Execute("30+19=49", $true, $true, $false, $false, $false, $true, 1, $false, "5+46=51", 2) | Out-Null
$d.Content.Find.Execute("99-62=37", $true, $true, $false, $false, $false, $true, 1, $false, "33-14=19", 2) | Out-Null
$d.Content.Find.Execute("31+62=93", $true, $true, $false, $false, $false, $true, 1, $false, "94-70=24", 2) | Out-Null
$d.Content.Find.Execute("92-73=19", $true, $true, $false, $false, $false, $true, 1, $false, "53+29=82", 2) | Out-Null
$d.Content.Find.Execute("95-72=23", $true, $true, $false, $false, $false, $true, 1, $false, "68-35=33", 2) | Out-Null
$d.Content.Find.Execute("0+41=41", $true, $true, $false, $false, $false, $true, 1, $false, "33+19=52", 2) | Out-Null
$d.Content.Find.Execute("12+30=42", $true, $true, $false, $false, $false, $true, 1, $false, "48-35=13", 2) | Out-Null
$d.Content.Find.Execute("63-45=18", $true, $true, $false, $false, $false, $true, 1, $false, "22+61=83", 2) | Out-Null
$d.Content.Find.Execute("26+39=65", $true, $true, $false, $false, $false, $true, 1, $false, "29-4=25", 2) | Out-Null
$d.Content.Find.Execute("26-18=8", $true, $true, $false, $false, $false, $true, 1, $false, "15+8=23", 2) | Out-Null
$d.Content.Find.Execute("23+53=76", $true, $true, $false, $false, $false, $true, 1, $false, "84-61=23", 2) | Out-Null
$d.Content.Find.Execute("99-96=3", $true, $true, $false, $false, $false, $true, 1, $false, "1+27=28", 2) | Out-Null
$d.Content.Find.Execute("60-40=20", $true, $true, $false, $false, $false, $true, 1, $false, "74-45=29", 2) | Out-Null
$d.Content.Find.Execute("49-43=6", $true, $true, $false, $false, $false, $true, 1, $false, "41+41=82", 2) | Out-Null
$d.Content.Find.Execute("97-68=29", $true, $true, $false, $false, $false, $true, 1, $false, "76-66=10", 2) | Out-Null
$d.Content.Find.Execute("85-53=32", $true, $true, $false, $false, $false, $true, 1, $false, "74+2=76", 2) | Out-Null
$d.Content.Find.Execute("0+49=49", $true, $true, $false, $false, $false, $true, 1, $false, "60-46=14", 2) | Out-Null
$d.Content.Find.Execute("76-40=36", $true, $true, $false, $false, $false, $true, 1, $false, "74+8=82", 2) | Out-Null
$d.Content.Find.Execute("35+49=84", $true, $true, $false, $false, $false, $true, 1, $false, "23+55=78", 2) | Out-Null
$d.Content.Find.Execute("45+41=86", $true, $true, $false, $false, $false, $true, 1, $false, "20-5=15", 2) | Out-Null
$d.Content.Find.Execute("86-77=9", $true, $true, $false, $false, $false, $true, 1, $false, "86+12=98", 2) | Out-Null
$d.Content.Find.Execute("72-38=34", $true, $true, $false, $false, $false, $true, 1, $false, "7+57=64", 2) | Out-Null
$d.Content.Find.Execute("5+2=7", $true, $true, $false, $false, $false, $true, 1, $false, "84-67=17", 2) | Out-Null
$d.Content.Find.Execute("69-26=43", $true, $true, $false, $false, $false, $true, 1, $false, "25+18=43", 2) | Out-Null
$d.Content.Find.Execute("18+72=90", $true, $true, $false, $false, $false, $true, 1, $false, "11+85=96", 2) | Out-Null
$d.Content.Find.Execute("52+5=57", $true, $true, $false, $false, $false, $true, 1, $false, "55-4=51", 2) | Out-Null
$d.Content.Find.Execute("4+77=81", $true, $true, $false, $false, $false, $true, 1, $false, "59+11=70", 2) | Out-Null
$d.Content.Find.Execute("3+91=94", $true, $true, $false, $false, $false, $true, 1, $false, "1+56=57", 2) | Out-Null
$d.Content.Find.Execute("45+23=68", $true, $true, $false, $false, $false, $true, 1, $false, "60-24=36", 2) | Out-Null
$d.Content.Find.Execute("71-55=16", $true, $true, $false, $false, $false, $true, 1, $false, "73-42=31", 2) | Out-Null
$d.Content.Find.Execute("81-58=23", $true, $true, $false, $false, $false, $true, 1, $false, "68-20=48", 2) | Out-Null
$d.Content.Find.Execute("80-36=44", $true, $true, $false, $false, $false, $true, 1, $false, "26+63=89", 2) | Out-Null
$d.Content.Find.Execute("70-46=24", $true, $true, $false, $false, $false, $true, 1, $false, "53-45=8", 2) | Out-Null
$d.Content.Find.Execute("95-39=56", $true, $true, $false, $false, $false, $true, 1, $false, "39+46=85", 2) | Out-Null
$d.Content.Find.Execute("27-27=0", $true, $true, $false, $false, $false, $true, 1, $false, "53+44=97", 2) | Out-Null
$d.Content.Find.Execute("3+77=80", $true, $true, $false, $false, $false, $true, 1, $false, "1+24=25", 2) | Out-Null
$d.Content.Find.Execute("41+1=42", $true, $true, $false, $false, $false, $true, 1, $false, "46-34=12", 2) | Out-Null
$d.Content.Find.Execute("58+32=90", $true, $true, $false, $false, $false, $true, 1, $false, "61+23=84", 2) | Out-Null
$d.Content.Find.Execute("85-79=6", $true, $true, $false, $false, $false, $true, 1, $false, "77-16=61", 2) | Out-Null
$d.Content.Find.Execute("95-37=58", $true, $true, $false, $false, $false, $true, 1, $false, "35+24=59", 2) | Out-Null
$d.Content.Find.Execute("88-54=34", $true, $true, $false, $false, $false, $true, 1, $false, "91-50=41", 2) | Out-Null
$d.Content.Find.Execute("54-8=46", $true, $true, $false, $false, $false, $true, 1, $false, "53+24=77", 2) | Out-Null
$d.Content.Find.Execute("28+69=97", $true, $true, $false, $false, $false, $true, 1, $false, "35+1=36", 2) | Out-Null
$d.Content.Find.Execute("70-7=63", $true, $true, $false, $false, $false, $true, 1, $false, "38+58=96", 2) | Out-Null
$d.Content.Find.Execute("29+29=58", $true, $true, $false, $false, $false, $true, 1, $false, "72+4=76", 2) | Out-Null
$d.Content.Find.Execute("65+19=84", $true, $true, $false, $false, $false, $true, 1, $false, "89-43=46", 2) | Out-Null
$d.Content.Find.Execute("45+35=80", $true, $true, $false, $false, $false, $true, 1, $false, "6+63=69", 2) | Out-Null
$d.Content.Find.Execute("12+57=69", $true, $true, $false, $false, $false, $true, 1, $false, "48-28=20", 2) | Out-Null
$d.Content.Find.Execute("23+61=84", $true, $true, $false, $false, $false, $true, 1, $false, "74-57=17", 2) | Out-Null
$d.Content.Find.Execute("15-5=10", $true, $true, $false, $false, $false, $true, 1, $false, "32+11=43", 2) | Out-Null
$d.Content.Find.Execute("25-18=7", $true, $true, $false, $false, $false, $true, 1, $false, "78-27=51", 2) | Out-Null
$d.Content.Find.Execute("20+75=95", $true, $true, $false, $false, $false, $true, 1, $false, "50-39=11", 2) | Out-Null
$d.Content.Find.Execute("30-29=1", $true, $true, $false, $false, $false, $true, 1, $false, "93-13=80", 2) | Out-Null
$d.Content.Find.Execute("84-63=21", $true, $true, $false, $false, $false, $true, 1, $false, "71-71=0", 2) | Out-Null
$d.Content.Find.Execute("96-21=75", $true, $true, $false, $false, $false, $true, 1, $false, "84-20=64", 2) | Out-Null
$d.Content.Find.Execute("5+31=36", $true, $true, $false, $false, $false, $true, 1, $false, "61+31=92", 2) | Out-Null
$d.Content.Find.Execute("23+45=68", $true, $true, $false, $false, $false, $true, 1, $false, "53+41=94", 2) | Out-Null
$d.Content.Find.Execute("10+22=32", $true, $true, $false, $false, $false, $true, 1, $false, "16+26=42", 2) | Out-Null
$d.Content.Find.Execute("18+28=46", $true, $true, $false, $false, $false, $true, 1, $false, "0+72=72", 2) | Out-Null
$d.Content.Find.Execute("91-26=65", $true, $true, $false, $false, $false, $true, 1, $false, "39+57=96", 2) | Out-Null
$d.Content.Find.Execute("75-33=42", $true, $true, $false, $false, $false, $true, 1, $false, "55-36=19", 2) | Out-Null
$d.Content.Find.Execute("69+21=90", $true, $true, $false, $false, $false, $true, 1, $false, "65-18=47", 2) | Out-Null
$d.Content.Find.Execute("5+21=26", $true, $true, $false, $false, $false, $true, 1, $false, "72-7=65", 2) | Out-Null
$d.Content.Find.Execute("52+41=93", $true, $true, $false, $false, $false, $true, 1, $false, "6+73=79", 2) | Out-Null
$d.Content.Find.Execute("64-29=35", $true, $true, $false, $false, $false, $true, 1, $false, "31+18=49", 2) | Out-Null
$d.Content.Find.Execute("51-47=4", $true, $true, $false, $false, $false, $true, 1, $false, "47+3=50", 2) | Out-Null
$d.Content.Find.Execute("64+5=69", $true, $true, $false, $false, $false, $true, 1, $false, "6+60=66", 2) | Out-Null
$d.Content.Find.Execute("83-61=22", $true, $true, $false, $false, $false, $true, 1, $false, "13+26=39", 2) | Out-Null
$d.Content.Find.Execute("46-23=23", $true, $true, $false, $false, $false, $true, 1, $false, "94-79=15", 2) | Out-Null
$d.Content.Find.Execute("27+68=95", $true, $true, $false, $false, $false, $true, 1, $false, "52-21=31", 2) | Out-Null
$d.Content.Find.Execute("76-68=8", $true, $true, $false, $false, $false, $true, 1, $false, "6+64=70", 2) | Out-Null
$d.Content.Find.Execute("38+47=85", $true, $true, $false, $false, $false, $true, 1, $false, "31+61=92", 2) | Out-Null
$d.Content.Find.Execute("94-92=2", $true, $true, $false, $false, $false, $true, 1, $false, "1+17=18", 2) | Out-Null
$d.Content.Find.Execute("71+21=92", $true, $true, $false, $false, $false, $true, 1, $false, "80-68=12", 2) | Out-Null
$d.Content.Find.Execute("10+82=92", $true, $true, $false, $false, $false, $true, 1, $false, "68-64=4", 2) | Out-Null
